# "Major reorganisation of data objects / Completed implementation of tagging"
#
# Appends the next block of error-code rows (10051-10065) to Sheet1, including
# three brand-new "tagging" messages (10060-10062) with their Message/Level
# ("Success") columns populated, and updates the active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 10051-10059: code only, no message/level yet -----------------
$codeOnly = 10051
for ($r = 53; $r -le 61; $r++) {
    $ws.Cells.Item($r, 1).Value = $codeOnly
    $codeOnly++
}

# --- Rows 10060-10062: new tagging feature messages ---------------------
$ws.Range("A62").Value = 10060
$ws.Range("B62").Value = "message_10060_new_tag_created"
$ws.Range("D62").Value = "Success"

$ws.Range("A63").Value = 10061
$ws.Range("B63").Value = "message_10061_tag_deleted"
$ws.Range("D63").Value = "Success"

$ws.Range("A64").Value = 10062
$ws.Range("B64").Value = "message_10062_tag_collection_updated"
$ws.Range("D64").Value = "Success"

# --- Rows 10063-10065: code only, reserved for upcoming work -------------
$codeOnly = 10063
for ($r = 65; $r -le 67; $r++) {
    $ws.Cells.Item($r, 1).Value = $codeOnly
    $codeOnly++
}

# --- Update the view: scroll down and select the newest message cell -----
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B64").Select()
